$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.26"
$ws.Range("E2").Value = "'-1.58%"
$ws.Range("D3").Value = "'31.62"
$ws.Range("E3").Value = "'0.30%"
$ws.Range("D4").Value = "'5.080"
$ws.Range("E4").Value = "'-1.37%"
$ws.Range("D5").Value = "'0.08107"
$ws.Range("E5").Value = "'9.35%"
$ws.Range("D6").Value = "'2.601"
$ws.Range("E6").Value = "'21.35%"
$ws.Range("D7").Value = "'7.809"
$ws.Range("E7").Value = "'-1.38%"
$ws.Range("D8").Value = "'3.827"
$ws.Range("E8").Value = "'1.50%"
$ws.Range("D9").Value = "'0.9258"
$ws.Range("E9").Value = "'-0.10%"
$ws.Range("D10").Value = "'0.1758"
$ws.Range("E10").Value = "'1.65%"
$ws.Range("E11").Value = "'-3.21%"
$ws.Range("D12").Value = "'0.08916"
$ws.Range("E12").Value = "'9.07%"
$ws.Range("D13").Value = "'0.03031"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("E14").Value = "'0.64%"
$ws.Range("D15").Value = "'0.001499"
$ws.Range("E15").Value = "'-0.19%"
$ws.Range("D16").Value = "'0.006010"
$ws.Range("E16").Value = "'-1.67%"
$ws.Range("D17").Value = "'3.555"
$ws.Range("E17").Value = "'2.60%"
$ws.Range("E18").Value = "'0.67%"
$ws.Range("D19").Value = "'0.3229"
$ws.Range("E19").Value = "'-0.60%"
$ws.Range("D20").Value = "'0.1318"
$ws.Range("E20").Value = "'-2.22%"
$ws.Range("D21").Value = "'3.975"
$ws.Range("E21").Value = "'-14.46%"
$ws.Range("E22").Value = "'4.20%"
$ws.Range("D23").Value = "'0.04599"
$ws.Range("E23").Value = "'-0.96%"
$ws.Range("E24").Value = "'1.53%"
$ws.Range("D25").Value = "'0.004440"
$ws.Range("E25").Value = "'-1.18%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-7.63%"
$ws.Range("D27").Value = "'0.0003409"
$ws.Range("E27").Value = "'82.12%"
$ws.Range("D39").Value = "'0.01771"
$ws.Range("E39").Value = "'1.52%"
$ws.Range("D40").Value = "'0.04508"
$ws.Range("E40").Value = "'-0.72%"
$ws.Range("D41").Value = "'0.006801"
$ws.Range("E41").Value = "'-4.13%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'-0.14%"
$ws.Range("D43").Value = "'0.002141"
$ws.Range("E43").Value = "'-2.96%"
$ws.Range("D44").Value = "'0.009829"
$ws.Range("E44").Value = "'-10.41%"
$ws.Range("D45").Value = "'0.00006458"
$ws.Range("E45").Value = "'2.93%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("D47").Value = "'0.008740"
$ws.Range("E47").Value = "'24.90%"
$ws.Range("E48").Value = "'-57.44%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.02%"

# Reset style on touched cells back to Normal so no stray quote-prefix formatting remains
$ws.Range("D2:E50").Style = "Normal"
